# Reorders the "Periodo Mora" / "Valor Mora" rows of the account statement
# table (rows 16-20) so the periods are listed in descending order, and
# updates the "Valor Mora" figures to match the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order for the Periodo Mora (column E) values, rows 16-20
$periodos = @("1903", "1902", "1901", "1812", "1811")

# New order for the Valor Mora (column F) values, rows 16-20
$valores = @(26041, 31249, 31249, 31249, 31249)

for ($i = 0; $i -lt 5; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
